$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry below mirrors one changed row from the cryptos list refresh:
# coin name (B) / link (C) only change where two rows swapped rank order,
# while Price (D) and Volume(1h) (E) are refreshed for every listed coin.
$updates = @(
    @{ Row = 2; D = "81.368.27"; E = "  +4.55%  " },
    @{ Row = 3; D = "3.176.15"; E = "  +0.22%  " },
    @{ Row = 4; E = "  +0.05%  " },
    @{ Row = 5; D = "207.87"; E = "  +2.41%  " },
    @{ Row = 6; D = "632.96"; E = "  +0.64%  " },
    @{ Row = 7; D = "0.293"; E = "  +29.01%  " },
    @{ Row = 8; E = "  -0.01%  " },
    @{ Row = 9; D = "0.591"; E = "  +3.14%  " },
    @{ Row = 10; D = "3.173.30"; E = "  +0.19%  " },
    @{ Row = 11; D = "0.588"; E = "  +3.47%  " },
    @{ Row = 12; E = "  +16.41%  " },
    @{ Row = 13; E = "  +2.06%  " },
    @{ Row = 14; D = "5.36"; E = "  -0.81%  " },
    @{ Row = 15; D = "3.760.78"; E = "  +0.39%  " },
    @{ Row = 16; D = "31.99"; E = "  +1.42%  " },
    @{ Row = 17; D = "81.418.80"; E = "  +4.69%  " },
    @{ Row = 18; D = "3.171.40"; E = "  +1.14%  " },
    @{ Row = 19; D = "3.24"; E = "  +14.27%  " },
    @{ Row = 20; D = "14.20"; E = "  -0.90%  " },
    @{ Row = 21; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "9.20"; E = "  -1.96%  " },
    @{ Row = 22; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "439.18"; E = "  +2.43%  " },
    @{ Row = 23; D = "5.17"; E = "  +6.74%  " },
    @{ Row = 24; E = "  +6.11%  " },
    @{ Row = 25; D = "5.17"; E = "  +9.92%  " },
    @{ Row = 26; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "11.21"; E = "  +3.37%  " },
    @{ Row = 27; B = "WrappedeETH"; C = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; D = "3.345.64"; E = "  +0.68%  " },
    @{ Row = 28; D = "76.94"; E = "  +1.06%  " },
    @{ Row = 29; D = "0.999"; E = "  -0.34%  " },
    @{ Row = 30; E = "  +9.82%  " },
    @{ Row = 31; D = "9.14"; E = "  +3.62%  " },
    @{ Row = 32; E = "  +0.46%  " },
    @{ Row = 33; D = "560.33"; E = "  +7.89%  " },
    @{ Row = 34; E = "  +2.24%  " },
    @{ Row = 35; E = "  +3.13%  " },
    @{ Row = 36; E = "  +11.56%  " },
    @{ Row = 37; D = "0.139"; E = "  +28.69%  " },
    @{ Row = 38; D = "23.16"; E = "  +3.37%  " },
    @{ Row = 39; E = "  -0.01%  " },
    @{ Row = 40; E = "  +4.42%  " },
    @{ Row = 41; D = "3.10"; E = "  +21.86%  " },
    @{ Row = 42; D = "5.95"; E = "  +10.49%  " },
    @{ Row = 43; B = "WhiteBITCoin"; C = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D = "20.76"; E = "  +3.47%  " },
    @{ Row = 44; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "2.03"; E = "  +15.09%  " },
    @{ Row = 45; D = "160.10"; E = "  -2.12%  " },
    @{ Row = 46; E = "  +0.00%  " },
    @{ Row = 47; D = "188.90"; E = "  -3.75%  " },
    @{ Row = 48; E = "  +4.34%  " },
    @{ Row = 49; D = "44.53"; E = "  +3.84%  " },
    @{ Row = 50; D = "0.783"; E = "  -1.96%  " },
    @{ Row = 51; D = "4.27"; E = "  +4.62%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) {
        $ws.Range("B$($u.Row)").Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Range("C$($u.Row)").Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Column D stores prices as plain text in the source data (e.g. "14.20",
        # "9.20"); force a text number format first so Excel's COM Value setter
        # does not reinterpret the string as a number and drop the trailing
        # zero / significant digits.
        $ws.Range("D$($u.Row)").NumberFormat = "@"
        $ws.Range("D$($u.Row)").Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
